$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "DATE_OF_DISABLEMENT" in L1
$ws.Range("L1").Value = "DATE_OF_DISABLEMENT"

# Add new date value in L2, formatted like the other date columns (A2, D2)
$ws.Range("L2").Value2 = $ws.Range("D2").Value2
$ws.Range("L2").NumberFormat = $ws.Range("D2").NumberFormat

# Update K2 value from 0.01 to 0
$ws.Range("K2").Value = 0

# Update selection to L2
$ws.Range("L2").Select()
